$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue "D2" "54.983.90"
Set-TextValue "E2" "  -5.05%  "
Set-TextValue "D3" "2.915.91"
Set-TextValue "E3" "  -7.90%  "
Set-TextValue "E4" "  +0.07%  "
Set-TextValue "D5" "476.42"
Set-TextValue "E5" "  -9.69%  "
Set-TextValue "D6" "128.46"
Set-TextValue "E6" "  -4.37%  "
Set-TextValue "E7" "  +0.00%  "
Set-TextValue "D8" "2.913.47"
Set-TextValue "E8" "  -7.96%  "
Set-TextValue "E9" "  -8.94%  "
Set-TextValue "D10" "6.78"
Set-TextValue "E10" "  -8.15%  "
Set-TextValue "E11" "  -11.56%  "
Set-TextValue "E12" "  -13.25%  "
Set-TextValue "E13" "  -1.84%  "
Set-TextValue "D14" "3.416.06"
Set-TextValue "E14" "  -7.71%  "
Set-TextValue "D15" "23.81"
Set-TextValue "E15" "  -8.20%  "
Set-TextValue "D16" "54.966.55"
Set-TextValue "E16" "  -4.97%  "
Set-TextValue "D17" "2.920.12"
Set-TextValue "E17" "  -7.50%  "
Set-TextValue "E18" "  -11.81%  "
Set-TextValue "D19" "5.48"
Set-TextValue "E19" "  -5.88%  "
Set-TextValue "D20" "11.68"
Set-TextValue "E20" "  -11.24%  "
Set-TextValue "D21" "7.27"
Set-TextValue "E21" "  -10.32%  "
Set-TextValue "D22" "305.32"
Set-TextValue "E22" "  -12.48%  "
Set-TextValue "E24" "  -12.20%  "
Set-TextValue "D25" "59.17"
Set-TextValue "E25" "  -15.17%  "
Set-TextValue "D26" "1.00"
Set-TextValue "E26" "  +0.16%  "
Set-TextValue "E27" "  -8.20%  "
Set-TextValue "E28" "  -0.12%  "
Set-TextValue "E29" "  -15.22%  "
Set-TextValue "D30" "6.30"
Set-TextValue "E30" "  -8.79%  "
Set-TextValue "D31" "6.37"
Set-TextValue "E31" "  -8.50%  "
Set-TextValue "D32" "1.13"
Set-TextValue "E32" "  -7.06%  "
Set-TextValue "E33" "  -12.61%  "
Set-TextValue "D34" "18.97"
Set-TextValue "E34" "  -12.81%  "
Set-TextValue "D35" "145.76"
Set-TextValue "E35" "  -8.93%  "
Set-TextValue "D36" "4.25"
Set-TextValue "E36" "  -14.09%  "
Set-TextValue "D37" "5.48"
Set-TextValue "E37" "  -12.80%  "
Set-TextValue "E38" "  -12.39%  "
Set-TextValue "D39" "23.20"
Set-TextValue "E39" "  -10.95%  "
Set-TextValue "E40" "  -10.00%  "
Set-TextValue "D41" "2.946.29"
Set-TextValue "E41" "  -7.60%  "
Set-TextValue "E42" "  +0.03%  "
Set-TextValue "D43" "35.71"
Set-TextValue "E43" "  -12.14%  "
Set-TextValue "D44" "0.972"
Set-TextValue "E44" "  -10.84%  "
Set-TextValue "D45" "0.618"
Set-TextValue "E45" "  -11.65%  "
Set-TextValue "E46" "  -8.56%  "
Set-TextValue "E47" "  -12.72%  "
Set-TextValue "D48" "2.067.72"
Set-TextValue "E48" "  -8.99%  "
Set-TextValue "E49" "  -12.34%  "
Set-TextValue "E50" "  -6.90%  "
Set-TextValue "D51" "18.23"
Set-TextValue "E51" "  -11.36%  "
